# Updated cryptos list with GitHub Actions
# Applies the refreshed price / 1h-volume-change figures to the coin
# ranking table on the active worksheet (columns D = Price, E = Volume(1h)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.675.70'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '3.483.88'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '591.22'
$ws.Range("E5").Value = '  +2.42%  '
$ws.Range("D6").Value = '168.38'
$ws.Range("E6").Value = '  -1.51%  '
$ws.Range("D7").Value = '0.607'
$ws.Range("E7").Value = '  -1.78%  '
$ws.Range("D8").Value = '3.480.44'
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '0.191'
$ws.Range("E10").Value = '  +1.14%  '
$ws.Range("D11").Value = '6.77'
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("D12").Value = '0.573'
$ws.Range("E12").Value = '  -4.19%  '
$ws.Range("D13").Value = '46.65'
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("D14").Value = '0.0000278'
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = '4.038.41'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").Value = '614.28'
$ws.Range("E16").Value = '  -10.21%  '
$ws.Range("D17").Value = '8.32'
$ws.Range("E17").Value = '  -5.06%  '
$ws.Range("D18").Value = '3.483.02'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").Value = '68.683.77'
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("E20").Value = '  -2.11%  '
$ws.Range("D21").Value = '17.22'
$ws.Range("E21").Value = '  -0.93%  '
$ws.Range("D22").Value = '11.14'
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("D23").Value = '0.872'
$ws.Range("E23").Value = '  -3.83%  '
$ws.Range("D24").Value = '15.76'
$ws.Range("E24").Value = '  -4.59%  '
$ws.Range("D25").Value = '95.73'
$ws.Range("E25").Value = '  -1.87%  '
$ws.Range("D26").Value = '3.78'
$ws.Range("E26").Value = '  -1.16%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").Value = '2.61'
$ws.Range("E28").Value = '  -1.89%  '
$ws.Range("D29").Value = '9.12'
$ws.Range("E29").Value = '  -3.04%  '
$ws.Range("D30").Value = '32.83'
$ws.Range("E30").Value = '  -1.09%  '
$ws.Range("D31").Value = '8.41'
$ws.Range("E31").Value = '  -4.69%  '
$ws.Range("E32").Value = '  -2.75%  '
$ws.Range("E33").Value = '  -2.25%  '
$ws.Range("D34").Value = '6.80'
$ws.Range("E34").Value = '  -6.10%  '
$ws.Range("D35").Value = '572.86'
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("D36").Value = '10.69'
$ws.Range("E36").Value = '  -1.20%  '
$ws.Range("D37").Value = '3.50'
$ws.Range("E37").Value = '  -3.92%  '
$ws.Range("D38").Value = '56.99'
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("E39").Value = '  -3.60%  '
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("E41").Value = '  -0.45%  '
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("D43").Value = '3.394.47'
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("D44").Value = '0.323'
$ws.Range("E44").Value = '  -3.99%  '
$ws.Range("D45").Value = '32.54'
$ws.Range("E45").Value = '  -1.97%  '
$ws.Range("D46").Value = '0.0₃0690'
$ws.Range("E46").Value = '  -1.61%  '
$ws.Range("D47").Value = '2.84'
$ws.Range("E47").Value = '  -1.50%  '
$ws.Range("D48").Value = '2.55'
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("E49").Value = '  -3.53%  '
$ws.Range("D50").Value = '132.31'
$ws.Range("E50").Value = '  -1.53%  '
$ws.Range("D51").Value = '5.54'
$ws.Range("E51").Value = '  +9.12%  '
